# update color pallete uses
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Florida Space Rangers" row's Primary Color changes from #FCEC00 to #D0D02B
$ws.Range("C5").Value = "#D0D02B"

# Update the active selection to match the author's final cursor position
$ws.Range("C5").Select()
